$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values: B2, D2 updated; C2, E2 cleared
$ws.Range("B2").Value = 23.254962237594334
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 35.204401339886694
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 20.682618561610106
$ws.Range("C3").Value = -6.5016201590062561
$ws.Range("D3").Value = 31.903132892840148
$ws.Range("E3").Value = -0.7865532533276669

# Update the visible selection to match the new authored range
$ws.Range("B1:E3").Select()
